# New weekly price observation ("Fruta / hortaliza, semanal") is inserted
# as a new row 35 on the sheet, pushing the existing rows 35-117 down to
# 36-118 (dimension grows from A1:R117 to A1:R118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 35 - this shifts rows 35:117 down to 36:118 and
# carries the D-column date style (s="2") into the new row automatically.
$ws.Rows("35:35").Insert()

# Populate the newly inserted row 35 with this week's record.
$ws.Range("A35").Value = 2
$ws.Range("B35").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C35").Value = "Coquimbo"
$ws.Range("D35").Value = 45014
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 100112030
$ws.Range("G35").Value = "Poroto granado"
$ws.Range("H35").Value = "Sin especificar"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 600
$ws.Range("K35").Value = 28000
$ws.Range("L35").Value = 29000
$ws.Range("M35").Value = 28500
$ws.Range("N35").Value = "$/malla 25 kilos"
$ws.Range("O35").Value = "Provincia de Limarí"
$ws.Range("P35").Value = 1140
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"
